# Scheduled-runner price/profit refresh for the Leve profit sheets.
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) for the rows whose market data changed, across all
# class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3903.3098
$ws.Range("I64").Value = 3737.2
$ws.Range("J64").Value = 4064.8057
$ws.Range("K64").Value = 3737.2
$ws.Range("L64").Value = 4064.8057
$ws.Range("M64").Value = -3489.2
$ws.Range("N64").Value = -4560.8057

$ws.Range("H67").Value = 3903.3098
$ws.Range("I67").Value = 3737.2
$ws.Range("J67").Value = 4064.8057
$ws.Range("K67").Value = 3737.2
$ws.Range("L67").Value = 4064.8057
$ws.Range("M67").Value = -2879.2
$ws.Range("N67").Value = -5780.8057

$ws.Range("H74").Value = 3472.6667
$ws.Range("I74").Value = 3464.4443
$ws.Range("J74").Value = 3476.1904
$ws.Range("K74").Value = 3464.4443
$ws.Range("L74").Value = 3476.1904
$ws.Range("M74").Value = -2528.4443
$ws.Range("N74").Value = -5348.190399999999

$ws.Range("H76").Value = 2986.5173
$ws.Range("I76").Value = 2662.0833
$ws.Range("K76").Value = 2662.0833
$ws.Range("M76").Value = -2347.0833

$ws.Range("H77").Value = 3472.6667
$ws.Range("I77").Value = 3464.4443
$ws.Range("J77").Value = 3476.1904
$ws.Range("K77").Value = 17322.2215
$ws.Range("L77").Value = 17380.952
$ws.Range("M77").Value = -12642.2215
$ws.Range("N77").Value = -26740.952

$ws.Range("H79").Value = 2986.5173
$ws.Range("I79").Value = 2662.0833
$ws.Range("K79").Value = 2662.0833
$ws.Range("M79").Value = -1570.0833

$ws.Range("H107").Value = 1659.1177
$ws.Range("I107").Value = 1457.5
$ws.Range("J107").Value = 2600
$ws.Range("K107").Value = 1457.5
$ws.Range("L107").Value = 2600
$ws.Range("M107").Value = 462.5
$ws.Range("N107").Value = -6440

$ws.Range("H129").Value = 751.9231
$ws.Range("I129").Value = 515.8333
$ws.Range("J129").Value = 954.2857
$ws.Range("K129").Value = 1547.4999
$ws.Range("L129").Value = 2862.8571
$ws.Range("M129").Value = 3452.5001
$ws.Range("N129").Value = -12862.8571

$ws.Range("H132").Value = 2754.074
$ws.Range("I132").Value = 2017.3191
$ws.Range("K132").Value = 6051.9573
$ws.Range("M132").Value = -3521.9573

$ws.Range("H138").Value = 2268.7883
$ws.Range("I138").Value = 1126.1818
$ws.Range("J138").Value = 2438.6353
$ws.Range("K138").Value = 3378.5454
$ws.Range("L138").Value = 7315.9059
$ws.Range("M138").Value = 1761.4546
$ws.Range("N138").Value = -17595.9059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11123183
$ws.Range("I32").Value = 15156933
$ws.Range("J32").Value = 30369.25
$ws.Range("K32").Value = 15156933
$ws.Range("L32").Value = 30369.25
$ws.Range("M32").Value = -15156646
$ws.Range("N32").Value = -30943.25

$ws.Range("H63").Value = 2144.4
$ws.Range("I63").Value = 2105.5
$ws.Range("J63").Value = 2300
$ws.Range("K63").Value = 2105.5
$ws.Range("L63").Value = 2300
$ws.Range("M63").Value = -1419.5
$ws.Range("N63").Value = -3672

$ws.Range("H66").Value = 2144.4
$ws.Range("I66").Value = 2105.5
$ws.Range("J66").Value = 2300
$ws.Range("K66").Value = 10527.5
$ws.Range("L66").Value = 11500
$ws.Range("M66").Value = -7095.5
$ws.Range("N66").Value = -18364

$ws.Range("H122").Value = 1628.9333
$ws.Range("I122").Value = 1003.7778
$ws.Range("K122").Value = 3011.3334
$ws.Range("M122").Value = -561.3334

$ws.Range("H132").Value = 1414.9773
$ws.Range("I132").Value = 800.2
$ws.Range("J132").Value = 2732.3572
$ws.Range("K132").Value = 2400.6
$ws.Range("L132").Value = 8197.071599999999
$ws.Range("M132").Value = 129.3999999999996
$ws.Range("N132").Value = -13257.0716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1361.8572
$ws.Range("I107").Value = 1260.5454
$ws.Range("J107").Value = 1733.3334
$ws.Range("K107").Value = 1260.5454
$ws.Range("L107").Value = 1733.3334
$ws.Range("M107").Value = 659.4546
$ws.Range("N107").Value = -5573.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2801.8333
$ws.Range("I16").Value = 1505.5
$ws.Range("K16").Value = 1505.5
$ws.Range("M16").Value = -1218.5

$ws.Range("H99").Value = 2079.889
$ws.Range("I99").Value = 2201.8462
$ws.Range("J99").Value = 1762.8
$ws.Range("K99").Value = 2201.8462
$ws.Range("L99").Value = 1762.8
$ws.Range("M99").Value = -703.8462
$ws.Range("N99").Value = -4758.8

$ws.Range("H107").Value = 885.2941
$ws.Range("I107").Value = 662.2
$ws.Range("J107").Value = 978.25
$ws.Range("K107").Value = 662.2
$ws.Range("L107").Value = 978.25
$ws.Range("M107").Value = 1257.8
$ws.Range("N107").Value = -4818.25

$ws.Range("H113").Value = 2801.8333
$ws.Range("I113").Value = 1505.5
$ws.Range("K113").Value = 1505.5
$ws.Range("M113").Value = 664.5

$ws.Range("H126").Value = 2079.889
$ws.Range("I126").Value = 2201.8462
$ws.Range("J126").Value = 1762.8
$ws.Range("K126").Value = 6605.5386
$ws.Range("L126").Value = 5288.4
$ws.Range("M126").Value = -4135.5386
$ws.Range("N126").Value = -10228.4

$ws.Range("H132").Value = 1608.375
$ws.Range("I132").Value = 1268.4839
$ws.Range("K132").Value = 3805.4517
$ws.Range("M132").Value = -1275.4517

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H64").Value = 4000
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 12000
$ws.Range("N64").Value = -12540

$ws.Range("H67").Value = 4000
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 12000
$ws.Range("N67").Value = -13872

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6926.263
$ws.Range("I70").Value = 7299.9414
$ws.Range("J70").Value = 3750
$ws.Range("K70").Value = 7299.9414
$ws.Range("L70").Value = 3750
$ws.Range("M70").Value = -7029.9414
$ws.Range("N70").Value = -4290

$ws.Range("H73").Value = 6926.263
$ws.Range("I73").Value = 7299.9414
$ws.Range("J73").Value = 3750
$ws.Range("K73").Value = 7299.9414
$ws.Range("L73").Value = 3750
$ws.Range("M73").Value = -6363.9414
$ws.Range("N73").Value = -5622

$ws.Range("H122").Value = 4349107
$ws.Range("I122").Value = 5264122.5
$ws.Range("J122").Value = 2781.75
$ws.Range("K122").Value = 15792367.5
$ws.Range("L122").Value = 8345.25
$ws.Range("M122").Value = -15789917.5
$ws.Range("N122").Value = -13245.25

$ws.Range("H126").Value = 8643630
$ws.Range("I126").Value = 6668383.5
$ws.Range("J126").Value = 18519862
$ws.Range("K126").Value = 20005150.5
$ws.Range("L126").Value = 55559586
$ws.Range("M126").Value = -20002680.5
$ws.Range("N126").Value = -55564526

$ws.Range("H132").Value = 2724.7568
$ws.Range("I132").Value = 2519.0417
$ws.Range("J132").Value = 3104.5386
$ws.Range("K132").Value = 7557.125100000001
$ws.Range("L132").Value = 9313.6158
$ws.Range("M132").Value = -5027.125100000001
$ws.Range("N132").Value = -14373.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1200.6957
$ws.Range("I61").Value = 1121.6316
$ws.Range("J61").Value = 1576.25
$ws.Range("K61").Value = 1121.6316
$ws.Range("L61").Value = 1576.25
$ws.Range("M61").Value = -919.6315999999999
$ws.Range("N61").Value = -1980.25

$ws.Range("H113").Value = 1200.6957
$ws.Range("I113").Value = 1121.6316
$ws.Range("J113").Value = 1576.25
$ws.Range("K113").Value = 1121.6316
$ws.Range("L113").Value = 1576.25
$ws.Range("M113").Value = 1048.3684
$ws.Range("N113").Value = -5916.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 22306.75
$ws.Range("J117").Value = 22306.75
$ws.Range("L117").Value = 22306.75
$ws.Range("N117").Value = -31484.75

$ws.Range("H132").Value = 2319.5
$ws.Range("I132").Value = 1408.8889
$ws.Range("J132").Value = 3064.5454
$ws.Range("K132").Value = 4226.6667
$ws.Range("L132").Value = 9193.636200000001
$ws.Range("M132").Value = -1696.6667
$ws.Range("N132").Value = -14253.6362
